$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''62.995.83'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '''  -0.55%  '
$ws.Range('E2').Style = "Normal"
$ws.Range('D3').Value = '''3.225.64'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '''  -0.57%  '
$ws.Range('E3').Style = "Normal"
$ws.Range('E4').Value = '''  +0.19%  '
$ws.Range('E4').Style = "Normal"
$ws.Range('D5').Value = '''527.23'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '''  +3.53%  '
$ws.Range('E5').Style = "Normal"
$ws.Range('D6').Value = '''171.48'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '''  -2.47%  '
$ws.Range('E6').Style = "Normal"
$ws.Range('E7').Value = '''  +1.63%  '
$ws.Range('E7').Style = "Normal"
$ws.Range('E8').Value = '''  -0.05%  '
$ws.Range('E8').Style = "Normal"
$ws.Range('D9').Value = '''3.223.06'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '''  -0.59%  '
$ws.Range('E9').Style = "Normal"
$ws.Range('D10').Value = '''0.602'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '''  -0.84%  '
$ws.Range('E10').Style = "Normal"
$ws.Range('D11').Value = '''52.81'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '''  -6.62%  '
$ws.Range('E11').Style = "Normal"
$ws.Range('E12').Value = '''  +2.60%  '
$ws.Range('E12').Style = "Normal"
$ws.Range('E13').Value = '''  +0.72%  '
$ws.Range('E13').Style = "Normal"
$ws.Range('E14').Value = '''  +1.30%  '
$ws.Range('E14').Style = "Normal"
$ws.Range('D15').Value = '''3.741.48'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '''  -0.20%  '
$ws.Range('E15').Style = "Normal"
$ws.Range('E16').Value = '''  -1.38%  '
$ws.Range('E16').Style = "Normal"
$ws.Range('D17').Value = '''3.229.56'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '''  +0.01%  '
$ws.Range('E17').Style = "Normal"
$ws.Range('D18').Value = '''62.869.79'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '''  -0.14%  '
$ws.Range('E18').Style = "Normal"
$ws.Range('D19').Value = '''17.14'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '''  +1.02%  '
$ws.Range('E19').Style = "Normal"
$ws.Range('D20').Value = '''11.02'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '''  +3.31%  '
$ws.Range('E20').Style = "Normal"
$ws.Range('D21').Value = '''0.963'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '''  +3.52%  '
$ws.Range('E21').Style = "Normal"
$ws.Range('D22').Value = '''365.06'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '''  -0.07%  '
$ws.Range('E22').Style = "Normal"
$ws.Range('E23').Value = '''  +4.28%  '
$ws.Range('E23').Style = "Normal"
$ws.Range('B24').Value = 'RenderToken'
$ws.Range('C24').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D24').Value = '''11.17'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '''  +3.24%  '
$ws.Range('E24').Style = "Normal"
$ws.Range('B25').Value = 'Litecoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D25').Value = '''80.79'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '''  +2.69%  '
$ws.Range('E25').Style = "Normal"
$ws.Range('D26').Value = '''4.00'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '''  +7.03%  '
$ws.Range('E26').Style = "Normal"
$ws.Range('D27').Value = '''6.06'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '''  +2.24%  '
$ws.Range('E27').Style = "Normal"
$ws.Range('E28').Value = '''  +1.01%  '
$ws.Range('E28').Style = "Normal"
$ws.Range('D29').Value = '''11.22'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '''  +1.07%  '
$ws.Range('E29').Style = "Normal"
$ws.Range('E30').Value = '''  +0.09%  '
$ws.Range('E30').Style = "Normal"
$ws.Range('D31').Value = '''28.35'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '''  +1.32%  '
$ws.Range('E31').Style = "Normal"
$ws.Range('D32').Value = '''631.18'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '''  -1.45%  '
$ws.Range('E32').Style = "Normal"
$ws.Range('D33').Value = '''6.41'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '''  -2.34%  '
$ws.Range('E33').Style = "Normal"
$ws.Range('D34').Value = '''11.17'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '''  +1.85%  '
$ws.Range('E34').Style = "Normal"
$ws.Range('D35').Value = '''0.105'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '''  +3.74%  '
$ws.Range('E35').Style = "Normal"
$ws.Range('E36').Value = '''  -3.94%  '
$ws.Range('E36').Style = "Normal"
$ws.Range('E37').Value = '''  -0.10%  '
$ws.Range('E37').Style = "Normal"
$ws.Range('D38').Value = '''36.37'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '''  +3.07%  '
$ws.Range('E38').Style = "Normal"
$ws.Range('E39').Value = '''  +0.79%  '
$ws.Range('E39').Style = "Normal"
$ws.Range('D40').Value = '''0.999'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '''  +0.47%  '
$ws.Range('E40').Style = "Normal"
$ws.Range('E41').Value = '''  +10.31%  '
$ws.Range('E41').Style = "Normal"
$ws.Range('D42').Value = '''0.124'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '''  +1.40%  '
$ws.Range('E42').Style = "Normal"
$ws.Range('B43').Value = 'Fetch.AI'
$ws.Range('C43').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D43').Value = '''2.55'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '''  +9.86%  '
$ws.Range('E43').Style = "Normal"
$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D44').Value = '''2.867.07'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '''  +0.91%  '
$ws.Range('E44').Style = "Normal"
$ws.Range('D45').Value = '''2.95'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '''  +7.05%  '
$ws.Range('E45').Style = "Normal"
$ws.Range('E46').Value = '''  +4.50%  '
$ws.Range('E46').Style = "Normal"
$ws.Range('E47').Value = '''  +6.75%  '
$ws.Range('E47').Style = "Normal"
$ws.Range('E48').Value = '''  +4.63%  '
$ws.Range('E48').Style = "Normal"
$ws.Range('E49').Value = '''  -1.42%  '
$ws.Range('E49').Style = "Normal"
$ws.Range('E50').Value = '''  +2.14%  '
$ws.Range('E50').Style = "Normal"
$ws.Range('D51').Value = '''133.62'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '''  +1.99%  '
$ws.Range('E51').Style = "Normal"
